$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing row 29 (dataset shrank by one entry in the right-hand table)
$ws.Rows.Item(29).Delete()

# Remove leftover left-hand-table data in row 7 (left table now only spans rows 1,3-6)
$ws.Range("A7:H7").Clear()

# Apply new anchor-word scores (left table) and keyword scores (right table)
$ws.Cells.Item(1,1).Value = "negative"
$ws.Cells.Item(1,10).Value = "positive"
$ws.Cells.Item(3,1).Value = "crude"
$ws.Cells.Item(3,2).Value = 0.8235294117647058
$ws.Cells.Item(3,3).Value = 28
$ws.Cells.Item(3,4).Value = 28
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = $false
$ws.Cells.Item(3,8).Value = 6
$ws.Cells.Item(3,10).Value = "love"
$ws.Cells.Item(3,11).Value = 0.9565217391304348
$ws.Cells.Item(3,12).Value = 44
$ws.Cells.Item(3,13).Value = 44
$ws.Cells.Item(3,14).Value = 1
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = $false
$ws.Cells.Item(3,17).Value = 2
$ws.Cells.Item(4,1).Value = "crisis"
$ws.Cells.Item(4,2).Value = 0.613013698630137
$ws.Cells.Item(4,3).Value = 179
$ws.Cells.Item(4,4).Value = 179
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = $false
$ws.Cells.Item(4,8).Value = 113
$ws.Cells.Item(4,10).Value = "best"
$ws.Cells.Item(4,11).Value = 0.9322033898305084
$ws.Cells.Item(4,12).Value = 55
$ws.Cells.Item(4,13).Value = 55
$ws.Cells.Item(4,14).Value = 1
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,16).Value = $false
$ws.Cells.Item(4,17).Value = 4
$ws.Cells.Item(5,1).Value = "panic"
$ws.Cells.Item(5,2).Value = 0.1976744186046512
$ws.Cells.Item(5,3).Value = 102
$ws.Cells.Item(5,4).Value = 102
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = $false
$ws.Cells.Item(5,8).Value = 414
$ws.Cells.Item(5,10).Value = "interesting"
$ws.Cells.Item(5,11).Value = 0.9090909090909091
$ws.Cells.Item(5,12).Value = 30
$ws.Cells.Item(5,13).Value = 30
$ws.Cells.Item(5,14).Value = 1
$ws.Cells.Item(5,15).Value = 0
$ws.Cells.Item(5,16).Value = $false
$ws.Cells.Item(5,17).Value = 3
$ws.Cells.Item(6,1).Value = "sc"
$ws.Cells.Item(6,2).Value = 0.1693121693121693
$ws.Cells.Item(6,3).Value = 32
$ws.Cells.Item(6,4).Value = 32
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = $false
$ws.Cells.Item(6,8).Value = 157
$ws.Cells.Item(6,10).Value = "great"
$ws.Cells.Item(6,11).Value = 0.875
$ws.Cells.Item(6,12).Value = 98
$ws.Cells.Item(6,13).Value = 98
$ws.Cells.Item(6,14).Value = 1
$ws.Cells.Item(6,15).Value = 0
$ws.Cells.Item(6,16).Value = $false
$ws.Cells.Item(6,17).Value = 14
$ws.Cells.Item(7,10).Value = "positive"
$ws.Cells.Item(7,11).Value = 0.7931034482758621
$ws.Cells.Item(7,12).Value = 46
$ws.Cells.Item(7,13).Value = 46
$ws.Cells.Item(7,14).Value = 1
$ws.Cells.Item(7,15).Value = 0
$ws.Cells.Item(7,16).Value = $false
$ws.Cells.Item(7,17).Value = 12
$ws.Cells.Item(8,10).Value = "thanks"
$ws.Cells.Item(8,11).Value = 0.7926829268292683
$ws.Cells.Item(8,12).Value = 65
$ws.Cells.Item(8,13).Value = 65
$ws.Cells.Item(8,14).Value = 1
$ws.Cells.Item(8,15).Value = 0
$ws.Cells.Item(8,16).Value = $false
$ws.Cells.Item(8,17).Value = 17
$ws.Cells.Item(9,10).Value = "free"
$ws.Cells.Item(9,11).Value = 0.7833333333333333
$ws.Cells.Item(9,12).Value = 94
$ws.Cells.Item(9,13).Value = 94
$ws.Cells.Item(9,14).Value = 1
$ws.Cells.Item(9,15).Value = 0
$ws.Cells.Item(9,16).Value = $false
$ws.Cells.Item(9,17).Value = 26
$ws.Cells.Item(10,10).Value = "thank"
$ws.Cells.Item(10,11).Value = 0.78125
$ws.Cells.Item(10,12).Value = 100
$ws.Cells.Item(10,13).Value = 100
$ws.Cells.Item(10,14).Value = 1
$ws.Cells.Item(10,15).Value = 0
$ws.Cells.Item(10,16).Value = $false
$ws.Cells.Item(10,17).Value = 28
$ws.Cells.Item(11,10).Value = "special"
$ws.Cells.Item(11,11).Value = 0.7777777777777778
$ws.Cells.Item(11,12).Value = 28
$ws.Cells.Item(11,13).Value = 28
$ws.Cells.Item(11,14).Value = 1
$ws.Cells.Item(11,15).Value = 0
$ws.Cells.Item(11,16).Value = $false
$ws.Cells.Item(11,17).Value = 8
$ws.Cells.Item(12,10).Value = "safe"
$ws.Cells.Item(12,11).Value = 0.7323943661971831
$ws.Cells.Item(12,12).Value = 104
$ws.Cells.Item(12,13).Value = 104
$ws.Cells.Item(12,14).Value = 1
$ws.Cells.Item(12,15).Value = 0
$ws.Cells.Item(12,16).Value = $false
$ws.Cells.Item(12,17).Value = 38
$ws.Cells.Item(13,10).Value = "confidence"
$ws.Cells.Item(13,11).Value = 0.7222222222222222
$ws.Cells.Item(13,12).Value = 26
$ws.Cells.Item(13,13).Value = 26
$ws.Cells.Item(13,14).Value = 1
$ws.Cells.Item(13,15).Value = 0
$ws.Cells.Item(13,16).Value = $false
$ws.Cells.Item(13,17).Value = 10
$ws.Cells.Item(14,10).Value = "good"
$ws.Cells.Item(14,11).Value = 0.70625
$ws.Cells.Item(14,12).Value = 113
$ws.Cells.Item(14,13).Value = 113
$ws.Cells.Item(14,14).Value = 1
$ws.Cells.Item(14,15).Value = 0
$ws.Cells.Item(14,16).Value = $false
$ws.Cells.Item(14,17).Value = 47
$ws.Cells.Item(15,10).Value = "support"
$ws.Cells.Item(15,11).Value = 0.6886792452830188
$ws.Cells.Item(15,12).Value = 73
$ws.Cells.Item(15,13).Value = 73
$ws.Cells.Item(15,14).Value = 1
$ws.Cells.Item(15,15).Value = 0
$ws.Cells.Item(15,16).Value = $false
$ws.Cells.Item(15,17).Value = 33
$ws.Cells.Item(16,10).Value = "safety"
$ws.Cells.Item(16,11).Value = 0.6862745098039216
$ws.Cells.Item(16,12).Value = 35
$ws.Cells.Item(16,13).Value = 35
$ws.Cells.Item(16,14).Value = 1
$ws.Cells.Item(16,15).Value = 0
$ws.Cells.Item(16,16).Value = $false
$ws.Cells.Item(16,17).Value = 16
$ws.Cells.Item(17,10).Value = "relief"
$ws.Cells.Item(17,11).Value = 0.62
$ws.Cells.Item(17,12).Value = 31
$ws.Cells.Item(17,13).Value = 31
$ws.Cells.Item(17,14).Value = 1
$ws.Cells.Item(17,15).Value = 0
$ws.Cells.Item(17,16).Value = $false
$ws.Cells.Item(17,17).Value = 19
$ws.Cells.Item(18,10).Value = "well"
$ws.Cells.Item(18,11).Value = 0.6170212765957447
$ws.Cells.Item(18,12).Value = 58
$ws.Cells.Item(18,13).Value = 58
$ws.Cells.Item(18,14).Value = 1
$ws.Cells.Item(18,15).Value = 0
$ws.Cells.Item(18,16).Value = $false
$ws.Cells.Item(18,17).Value = 36
$ws.Cells.Item(19,10).Value = "better"
$ws.Cells.Item(19,11).Value = 0.5873015873015873
$ws.Cells.Item(19,12).Value = 37
$ws.Cells.Item(19,13).Value = 37
$ws.Cells.Item(19,14).Value = 1
$ws.Cells.Item(19,15).Value = 0
$ws.Cells.Item(19,16).Value = $false
$ws.Cells.Item(19,17).Value = 26
$ws.Cells.Item(20,10).Value = "fresh"
$ws.Cells.Item(20,11).Value = 0.5625
$ws.Cells.Item(20,12).Value = 27
$ws.Cells.Item(20,13).Value = 27
$ws.Cells.Item(20,14).Value = 1
$ws.Cells.Item(20,15).Value = 0
$ws.Cells.Item(20,16).Value = $false
$ws.Cells.Item(20,17).Value = 21
$ws.Cells.Item(21,10).Value = "heroes"
$ws.Cells.Item(21,11).Value = 0.5319148936170213
$ws.Cells.Item(21,12).Value = 25
$ws.Cells.Item(21,13).Value = 25
$ws.Cells.Item(21,14).Value = 1
$ws.Cells.Item(21,15).Value = 0
$ws.Cells.Item(21,16).Value = $false
$ws.Cells.Item(21,17).Value = 22
$ws.Cells.Item(22,10).Value = "hand"
$ws.Cells.Item(22,11).Value = 0.5091383812010444
$ws.Cells.Item(22,12).Value = 195
$ws.Cells.Item(22,13).Value = 195
$ws.Cells.Item(22,14).Value = 1
$ws.Cells.Item(22,15).Value = 0
$ws.Cells.Item(22,16).Value = $false
$ws.Cells.Item(22,17).Value = 188
$ws.Cells.Item(23,10).Value = "like"
$ws.Cells.Item(23,11).Value = 0.4588235294117647
$ws.Cells.Item(23,12).Value = 156
$ws.Cells.Item(23,13).Value = 156
$ws.Cells.Item(23,14).Value = 1
$ws.Cells.Item(23,15).Value = 0
$ws.Cells.Item(23,16).Value = $false
$ws.Cells.Item(23,17).Value = 184
$ws.Cells.Item(24,10).Value = "care"
$ws.Cells.Item(24,11).Value = 0.449438202247191
$ws.Cells.Item(24,12).Value = 40
$ws.Cells.Item(24,13).Value = 40
$ws.Cells.Item(24,14).Value = 1
$ws.Cells.Item(24,15).Value = 0
$ws.Cells.Item(24,16).Value = $false
$ws.Cells.Item(24,17).Value = 49
$ws.Cells.Item(25,10).Value = "help"
$ws.Cells.Item(25,11).Value = 0.4135593220338983
$ws.Cells.Item(25,12).Value = 122
$ws.Cells.Item(25,13).Value = 122
$ws.Cells.Item(25,14).Value = 1
$ws.Cells.Item(25,15).Value = 0
$ws.Cells.Item(25,16).Value = $false
$ws.Cells.Item(25,17).Value = 173
$ws.Cells.Item(26,10).Value = "protect"
$ws.Cells.Item(26,11).Value = 0.3561643835616438
$ws.Cells.Item(26,12).Value = 26
$ws.Cells.Item(26,13).Value = 26
$ws.Cells.Item(26,14).Value = 1
$ws.Cells.Item(26,15).Value = 0
$ws.Cells.Item(26,16).Value = $false
$ws.Cells.Item(26,17).Value = 47
$ws.Cells.Item(27,10).Value = "increase"
$ws.Cells.Item(27,11).Value = 0.3461538461538461
$ws.Cells.Item(27,12).Value = 27
$ws.Cells.Item(27,13).Value = 27
$ws.Cells.Item(27,14).Value = 1
$ws.Cells.Item(27,15).Value = 0
$ws.Cells.Item(27,16).Value = $false
$ws.Cells.Item(27,17).Value = 51
$ws.Cells.Item(28,10).Value = "please"
$ws.Cells.Item(28,11).Value = 0.2928870292887029
$ws.Cells.Item(28,12).Value = 70
$ws.Cells.Item(28,13).Value = 70
$ws.Cells.Item(28,14).Value = 1
$ws.Cells.Item(28,15).Value = 0
$ws.Cells.Item(28,16).Value = $false
$ws.Cells.Item(28,17).Value = 169
